$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.959.93"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.842.06"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").Value = "'1.015"
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "'308.79"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "'0.4754"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("D8").Value = "'0.3670"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "'0.07206"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "'0.9298"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "'19.80"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "'0.07739"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.819.84"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").Value = "'5.375"
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "'6.458"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "'88.81"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "'1.017"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "'0.000008649"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "27.017.06"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'14.52"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "'5.051"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'1.927"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "'18.22"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").Value = "'1.994"
$ws.Range("E27").Value = "  -3.20%  "
$ws.Range("D28").Value = "'114.40"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").Value = "'4.960"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'3.308"
$ws.Range("E31").Value = "  +3.98%  "
$ws.Range("D32").Value = "'1.176"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.497"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7356"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").Value = "'2.680"
$ws.Range("E35").Value = "  -6.15%  "
$ws.Range("D36").Value = "'1.109"
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("D37").Value = "'0.01968"
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").Value = "'0.05257"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("D39").Value = "'2.975"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'0.5252"
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("D41").Value = "'7.015"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D43").Value = "'8.269"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("D44").Value = "'10.56"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").Value = "'0.4729"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").Value = "'101.65"
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("D48").Value = "'1.607"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "'65.59"
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").Value = "'0.06061"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'0.8913"
$ws.Range("E51").Value = "  +3.25%  "
